# Apply the "Plague Doctor" full-description row insertion to the "Text" sheet.
#
# Summary of the change:
#   - A new row is inserted before the existing row 549 (i.e. it becomes the
#     new row 548), which pushes every row from the old 549 onward down by
#     one. This matches the commit's intent of cleanly inserting the new
#     localization entry rather than overwriting/duplicating existing rows.
#   - The new row 548 gets three cells populated:
#       A548 = "plagueDoctorFullDesc"                (key)
#       B548 = English description (wrapped text)
#       M548 = Japanese description (wrapped text)
#   - Row 548 height is set to 120 and both text cells use wrap-text
#     formatting (matching the style already used by similar multi-line
#     description rows elsewhere in the sheet).
#   - The active selection is moved to A548 to reflect where the author was
#     last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Text")

# Insert a new row at position 549; existing row 549 (and everything below)
# shifts down to 550, etc. The newly-inserted, now-empty row becomes 548.
$ws.Rows.Item(549).Insert()

# Key name for this new localization entry.
$ws.Range("A548").Value = "plagueDoctorFullDesc"

# Populate Japanese (column M) before English (column B) so that the shared
# string table receives the Japanese text first, matching the original
# authoring order.
$ws.Range("M548").Value = "・誰か1人に感染病を付与できる`n・感染した人が一定時間の間非感染者の近くにいるとその人も感染者となる`n・感染状況は備蓄される（会議や離れても時間のリセットがされない）`n・ペスト医師がキルされたら、キルした人は感染者となる`n・生存者全員が感染者となることが勝利条件"

$ws.Range("B548").Value = "The Plague Doctor is a neutral role whose goal is to infect every living player.`nThey start by choosing one player to infect, after which anyone who spends a set`namount of time in range of the infected player becomes infected themselves.`nInfection progress is cumulative, and does not reset with distance or after meetings."

# Match the wrap-text formatting used by the other long description rows.
$ws.Range("B548").WrapText = $true
$ws.Range("M548").WrapText = $true

# Row needs to be tall enough to show the wrapped, multi-line text.
$ws.Rows.Item(548).RowHeight = 120

# Reflect where the author's selection ended up after adding the row.
$ws.Activate()
$ws.Range("A548").Select()
